$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values (Natmi re-run following Dr Hou advice)
$ws.Range("B2").Value = "Ccl21b"
$ws.Range("C2").Value = "Ackr2"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3964663333333333
$ws.Range("H2").Value = 1.189399
$ws.Range("I2").Value = 0.8220900069740165
$ws.Range("J2").Value = 0.8220900069740164
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 11.51723133333333
$ws.Range("N2").Value = 34.551694
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 4.566194476878444
$ws.Range("R2").Value = 41.09575029190599
$ws.Range("S2").Value = 0.8220900069740165
$ws.Range("T2").Value = 0.8220900069740164

# Add new row 3
$ws.Range("A3").Value = "sCs"
$ws.Range("B3").Value = "Ccl21b"
$ws.Range("C3").Value = "Ackr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.0858
$ws.Range("H3").Value = 0.2574
$ws.Range("I3").Value = 0.1779099930259836
$ws.Range("J3").Value = 0.1779099930259836
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 11.51723133333333
$ws.Range("N3").Value = 34.551694
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.9881784484
$ws.Range("R3").Value = 8.8936060356
$ws.Range("S3").Value = 0.1779099930259836
$ws.Range("T3").Value = 0.1779099930259836
